# Apply edits described by the diff:
# 1. Remove the "February" header/week-range cells (K2 and J3:M3) by
#    deleting the now-unused columns J:M outright (they shift nothing in,
#    since they are the trailing columns) - this also shrinks the
#    "2024"/"January" merges (F1:M1 -> F1:I1, F2:J2 -> F2:I2) and the
#    sheet dimension (B1:M9 -> B1:I9) automatically.
# 2. Fix the Task 2.1 / Task 2.2 Start/End date text values, which were
#    still pointing at February while their highlighted week cells
#    (H8, I9) actually line up with the January week columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-unused February columns (J:M) ---
$ws.Range("J1:M1").EntireColumn.Delete()

# --- Fix Task 2.1 / Task 2.2 dates (Feb -> Jan) ---
$ws.Range("D8").Value = "01/15"
$ws.Range("E8").Value = "01/21"
$ws.Range("D9").Value = "01/22"
$ws.Range("E9").Value = "01/28"
